$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# ---------------------------------------------------------------------------
# Column A width (wide label column)
# ---------------------------------------------------------------------------
$ws.Columns("A:A").ColumnWidth = 76.9986979166667

# ---------------------------------------------------------------------------
# Row 1 - header row
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Settings/Patients"
$ws.Range("A1").Font.Size = 18

$ws.Range("B1").Value = 231
$ws.Range("C1").Value = "222_1"
$ws.Range("D1").Value = "222_2"
$ws.Range("E1").Value = 241
$ws.Range("B1:E1").Font.Size = 24
$ws.Range("B1:E1").HorizontalAlignment = $xlCenter

$ws.Range("L1").Value = 229
$ws.Range("L1").Font.Size = 24
$ws.Range("L1").HorizontalAlignment = $xlCenter

# ---------------------------------------------------------------------------
# Row 2 - "124 with delta" (black band)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "124 with delta"
$ws.Range("A2").Font.Size = 18
$ws.Range("A2").Interior.ThemeColor = 1

$ws.Range("B2").Value = 88.77
$ws.Range("C2").Value = 75.34
$ws.Range("D2").Value = 70.15
$ws.Range("E2").Value = 87.87
$ws.Range("B2:E2").Font.Size = 24
$ws.Range("B2:E2").HorizontalAlignment = $xlCenter
$ws.Range("B2:E2").Interior.ThemeColor = 1

$ws.Range("L2").Value = 83.69
$ws.Range("L2").Font.Size = 24
$ws.Range("L2").HorizontalAlignment = $xlCenter
$ws.Range("L2").Interior.ThemeColor = 1

# ---------------------------------------------------------------------------
# Row 3 - "124 without delta" (black band)
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "124 without delta"
$ws.Range("A3").Font.Size = 18
$ws.Range("A3").Interior.ThemeColor = 1

$ws.Range("B3").Value = 89.14
$ws.Range("C3").Value = 72.85
$ws.Range("D3").Value = 70.44
$ws.Range("B3:D3").Font.Size = 24
$ws.Range("B3:D3").HorizontalAlignment = $xlCenter
$ws.Range("B3:D3").Interior.ThemeColor = 1

$ws.Range("E3").Clear()

$ws.Range("L3").Value = 83.74
$ws.Range("L3").Font.Size = 24
$ws.Range("L3").HorizontalAlignment = $xlCenter
$ws.Range("L3").Interior.ThemeColor = 1

# ---------------------------------------------------------------------------
# Row 4 - "Original monthly epoch, original Python split (90 without delta)" (white band)
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "Original monthly epoch, original Python split (90 without delta)"
$ws.Range("A4").Font.Size = 18
$ws.Range("A4").Interior.ThemeColor = 2

$ws.Range("B4").Value = 88.54
$ws.Range("C4").Value = 73.47
$ws.Range("D4").Value = 70.26
$ws.Range("E4").Value = 87.87
$ws.Range("B4:E4").Font.Size = 24
$ws.Range("B4:E4").HorizontalAlignment = $xlCenter
$ws.Range("B4:E4").Interior.ThemeColor = 2

$ws.Range("L4").Value = 83.12
$ws.Range("L4").Font.Size = 24
$ws.Range("L4").HorizontalAlignment = $xlCenter
$ws.Range("L4").Interior.ThemeColor = 2

# ---------------------------------------------------------------------------
# Row 5 - "Original monthly epoch, NEW even split"
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "Original monthly epoch, NEW even split"
$ws.Range("A5").Font.Size = 18

$ws.Range("B5").Value = 83.5
$ws.Range("C5").Value = 65.9
$ws.Range("D5").Value = 77.25
$ws.Range("E5").Value = 87.84
$ws.Range("B5:E5").Font.Size = 24
$ws.Range("B5:E5").HorizontalAlignment = $xlCenter

$ws.Range("F5:G5").Font.Size = 24
$ws.Range("F5:G5").HorizontalAlignment = $xlCenter

# ---------------------------------------------------------------------------
# Row 6 - blank filler row (same per-column formatting as row 5)
# ---------------------------------------------------------------------------
$ws.Rows(6).RowHeight = 31
$ws.Range("B6:G6").Font.Size = 24
$ws.Range("B6:G6").HorizontalAlignment = $xlCenter

# ---------------------------------------------------------------------------
# Row 7 - "NEW weekly epoch, NEW even split"
# ---------------------------------------------------------------------------
$ws.Rows(7).RowHeight = 31
$ws.Range("A7").Value = "NEW weekly epoch, NEW even split"
$ws.Range("A7").Font.Size = 18

$ws.Range("B7").Value = 60.94
$ws.Range("C7").Value = 73.78
$ws.Range("D7").Value = 59.3
$ws.Range("E7").Value = 72.22
$ws.Range("B7:E7").Font.Size = 24
$ws.Range("B7:E7").HorizontalAlignment = $xlCenter

$ws.Range("F7:G7").Font.Size = 24
$ws.Range("F7:G7").HorizontalAlignment = $xlCenter

# ---------------------------------------------------------------------------
# Row 8 - "NEW biweekly epoch, NEW even split"
# ---------------------------------------------------------------------------
$ws.Rows(8).RowHeight = 31
$ws.Range("A8").Value = "NEW biweekly epoch, NEW even split"
$ws.Range("A8").Font.Size = 18

$ws.Range("B8").Value = 66.58
$ws.Range("C8").Value = 72.73
$ws.Range("D8").Value = 70.25
$ws.Range("E8").Value = 87.74
$ws.Range("B8:E8").Font.Size = 24
$ws.Range("B8:E8").HorizontalAlignment = $xlCenter

$ws.Range("F8:G8").Font.Size = 24
$ws.Range("F8:G8").HorizontalAlignment = $xlCenter

# ---------------------------------------------------------------------------
# Rows 9-13 - blank filler rows (same per-column formatting as row 5)
# ---------------------------------------------------------------------------
foreach ($rr in 9..13) {
    $ws.Rows($rr).RowHeight = 31
    $rng = $ws.Range("B${rr}:G${rr}")
    $rng.Font.Size = 24
    $rng.HorizontalAlignment = $xlCenter
}

# ---------------------------------------------------------------------------
# Selection / view
# ---------------------------------------------------------------------------
$ws.Range("A5").Select()
